# ---------------------------------------------------------------------------
# Commit: "Update data: add rebounds & 3-pointers"
#
# Adds two new per-game stat sheets ("Rebounds", "3PM") right after the
# existing "Assists" sheet, and two new per-player average sheets
# ("Avg Rebounds", "Avg 3PM") right after the existing "Avg Assists" sheet,
# populating each with the season-to-date data.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Fill-Matrix($ws, $topLeftCell, $matrix) {
    $nRows = $matrix.Count
    $nCols = $matrix[0].Count
    $arr = New-Object 'object[,]' $nRows,$nCols
    for ($r = 0; $r -lt $nRows; $r++) {
        for ($c = 0; $c -lt $nCols; $c++) {
            $arr[$r,$c] = $matrix[$r][$c]
        }
    }
    $rng = $ws.Range($topLeftCell).Resize($nRows, $nCols)
    $rng.Value = $arr
    return $rng
}

function Style-Header($ws, $headerRange) {
    $srcHeader = $wb.Worksheets.Item('Points').Range('A1')
    $srcHeader.Copy()
    $headerRange.PasteSpecial(-4122, [System.Reflection.Missing]::Value, $false, $false)
    $ws.Application.CutCopyMode = $false
}

# --- Insert 'Rebounds' and '3PM' sheets right after 'Assists' ---
$afterAssists = $wb.Worksheets.Item('Assists')
$wsRebounds = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterAssists)
$wsRebounds.Name = 'Rebounds'

$afterRebounds = $wb.Worksheets.Item('Rebounds')
$ws3PM = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterRebounds)
$ws3PM.Name = '3PM'

# --- Insert 'Avg Rebounds' and 'Avg 3PM' sheets right after 'Avg Assists' ---
$afterAvgAssists = $wb.Worksheets.Item('Avg Assists')
$wsAvgRebounds = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterAvgAssists)
$wsAvgRebounds.Name = 'Avg Rebounds'

$afterAvgRebounds = $wb.Worksheets.Item('Avg Rebounds')
$wsAvg3PM = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterAvgRebounds)
$wsAvg3PM.Name = 'Avg 3PM'

# --- Populate 'Rebounds' ---
$reboundsData = @(
        @('Game Time (PST)','Opponent','Gradey Dick','Jonathan Mogbo','Brandon Ingram','Scottie Barnes','Immanuel Quickley','RJ Barrett','Collin Murray-Boyles','Ja''Kobe Walter','Jakob Poeltl','Jamal Shead','Ochai Agbaji','Sandro Mamukelashvili','Jamison Battle'),
        @('''2025-10-22','ATL',5,4,9,6,4,8,0,0,6,2,3,4,3),
        @('''2025-10-24','MIL',4,2,6,5,7,6,0,0,5,0,1,3,2),
        @('''2025-10-26','DAL',2,0,6,11,5,4,2,0,8,0,3,1,1),
        @('''2025-10-27','SAS',0,0,3,3,3,2,3,1,2,2,1,0,0),
        @('''2025-10-29','HOU',3,1,2,5,3,2,1,2,0,2,0,1,0),
        @('''2025-10-31','CLE',0,0,8,10,1,3,4,1,0,7,7,6,1),
        @('''2025-11-02','MEM',4,0,7,12,4,6,9,2,0,0,2,3,0),
        @('''2025-11-04','MIL',1,2,8,3,6,8,1,3,9,2,0,7,0),
        @('''2025-11-07','ATL',1,0,6,10,5,5,4,0,10,1,0,6,0),
        @('''2025-11-08','PHI',1,0,8,5,6,3,7,0,0,3,2,3,0)
    )
$rngRebounds = Fill-Matrix $wsRebounds 'A1' $reboundsData
Style-Header $wsRebounds ($wsRebounds.Range('A1:O1'))

# --- Populate '3PM' ---
$threePmData = @(
        @('Game Time (PST)','Opponent','Gradey Dick','Jonathan Mogbo','Brandon Ingram','Scottie Barnes','Immanuel Quickley','RJ Barrett','Collin Murray-Boyles','Ja''Kobe Walter','Jakob Poeltl','Jamal Shead','Ochai Agbaji','Sandro Mamukelashvili','Jamison Battle'),
        @('''2025-10-22','ATL',2,0,0,0,0,2,0,0,0,1,0,1,0),
        @('''2025-10-24','MIL',1,0,3,2,1,1,0,0,0,0,0,0,3),
        @('''2025-10-26','DAL',2,0,2,3,1,1,0,0,0,1,1,1,1),
        @('''2025-10-27','SAS',0,0,0,2,2,2,3,1,0,1,0,1,0),
        @('''2025-10-29','HOU',0,0,5,4,4,2,3,0,0,2,0,1,0),
        @('''2025-10-31','CLE',0,0,0,1,0,3,0,1,0,1,0,1,6),
        @('''2025-11-02','MEM',1,0,1,3,2,2,0,0,0,2,0,0,0),
        @('''2025-11-04','MIL',3,0,0,3,2,4,0,1,0,1,0,3,0),
        @('''2025-11-07','ATL',3,0,2,0,3,2,1,0,0,0,0,1,0),
        @('''2025-11-08','PHI',1,0,0,2,5,2,1,1,0,3,0,1,0)
    )
$rng3PM = Fill-Matrix $ws3PM 'A1' $threePmData
Style-Header $ws3PM ($ws3PM.Range('A1:O1'))

# --- Populate 'Avg Rebounds' ---
$avgReboundsData = @(
        @('Player','Avg Rebounds'),
        @('Scottie Barnes',7),
        @('Jakob Poeltl',6.666666666666667),
        @('Brandon Ingram',6.3),
        @('RJ Barrett',4.7),
        @('Immanuel Quickley',4.4),
        @('Collin Murray-Boyles',3.875),
        @('Sandro Mamukelashvili',3.4),
        @('Ochai Agbaji',2.111111111111111),
        @('Gradey Dick',2.1),
        @('Jamal Shead',1.9),
        @('Jonathan Mogbo',1.5),
        @('Ja''Kobe Walter',1.285714285714286),
        @('Jamison Battle',0.875)
    )
$rngAvgRebounds = Fill-Matrix $wsAvgRebounds 'A1' $avgReboundsData
Style-Header $wsAvgRebounds ($wsAvgRebounds.Range('A1:B1'))

# --- Populate 'Avg 3PM' ---
$avg3PMData = @(
        @('Player','Avg 3PM'),
        @('RJ Barrett',2.1),
        @('Scottie Barnes',2),
        @('Immanuel Quickley',2),
        @('Gradey Dick',1.3),
        @('Brandon Ingram',1.3),
        @('Jamison Battle',1.25),
        @('Jamal Shead',1.2),
        @('Collin Murray-Boyles',1),
        @('Sandro Mamukelashvili',1),
        @('Ja''Kobe Walter',0.5714285714285714),
        @('Ochai Agbaji',0.1111111111111111),
        @('Jonathan Mogbo',0),
        @('Jakob Poeltl',0)
    )
$rngAvg3PM = Fill-Matrix $wsAvg3PM 'A1' $avg3PMData
Style-Header $wsAvg3PM ($wsAvg3PM.Range('A1:B1'))

$wsRebounds.Select()
$wsRebounds.Range('A1').Select()

